# Implementing top panel buttons
# Insert 4 new rows (33-36) into the localization sheet for the new
# WARNING / approve / disapprove / exit-confirmation strings, pushing the
# existing taskName/taskContent rows down from 33-42 to 37-46.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before the current row 33 (taskName_3_1 row)
$ws.Rows("33:36").Insert()

# Copy the wrap-text formatting (style index 1) from the row above down
# onto the new B/C cells, matching the rest of the English/Polish columns.
$ws.Range("B32:C32").Copy()
$ws.Range("B33:C36").PasteSpecial(-4122)

# Row 33 - warning label
$ws.Range("B33").Value = "WARNING"
$ws.Range("A33").Value = "warningLabel_key"
$ws.Range("C33").Value = "OSTRZEŻENIE"

# Row 34 - approve button key
$ws.Range("A34").Value = "aprove_btn_key"

# Row 35 - disapprove button key
$ws.Range("A35").Value = "disaporve_btn_key"
$ws.Range("B35").Value = "No"

# Row 34 continued - Yes / Tak
$ws.Range("B34").Value = "Yes"
$ws.Range("C34").Value = "Tak"

# Row 35 continued - Nie
$ws.Range("C35").Value = "Nie"

# Row 36 - exit confirmation message
$ws.Range("A36").Value = "exitSystemMEssage_key"
$ws.Range("B36").Value = "Are you sure you want to exit HackSafe system?"
$ws.Range("C36").Value = "Czy na pewno chcesz opuścić HackSafe system?"

# Row 36 wraps to two lines in the real workbook (long Polish sentence).
$ws.Rows(36).RowHeight = 30

# Update the sheet view: pane/freeze still only splits off row 1, but the
# scrolled viewport + active selection moved.
$ws.Range("D31").Select()
